$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46029
$ws.Range("B2").Value = 14362.6898530377
$ws.Range("C2").Value = 14000.329322794
$ws.Range("D2").Value = 23595.86
$ws.Range("E2").Value = 9892.08898272747
$ws.Range("F2").Value = 12.3565960633957

$ws.Range("A3").Value = 46030
$ws.Range("B3").Value = 14318.3427522116
$ws.Range("C3").Value = 14059.6785854302
$ws.Range("D3").Value = 12075.86
$ws.Range("E3").Value = 9909.58153203659
$ws.Range("F3").Value = 495.558338227782

$ws.Range("A4").Value = 46031
$ws.Range("B4").Value = 13986.4405513841
$ws.Range("C4").Value = 12403.9093929529
$ws.Range("D4").Value = 12075.86
$ws.Range("E4").Value = 9665.60811131117
$ws.Range("F4").Value = 416.402396011001

$ws.Range("A5").Value = 46032
$ws.Range("B5").Value = 5646.03125881991
$ws.Range("C5").Value = 8755.1764392419
$ws.Range("D5").Value = 12075.86
$ws.Range("E5").Value = 8937.10455180606
$ws.Range("F5").Value = 234.017541293665

$ws.Range("A6").Value = 46033
$ws.Range("B6").Value = 5199.64162344583
$ws.Range("C6").Value = 8904.85436204661
$ws.Range("D6").Value = 12075.86
$ws.Range("E6").Value = 8617.63118729201
$ws.Range("F6").Value = 226.942731222443

$ws.Range("A7").Value = 46034
$ws.Range("B7").Value = 13378.423123507
$ws.Range("C7").Value = 13452.9930392632
$ws.Range("D7").Value = 12075.86
$ws.Range("E7").Value = 9107.12339392322
$ws.Range("F7").Value = 436.844018049432

$ws.Range("A8").Value = 46035
$ws.Range("B8").Value = 13378.423123507
$ws.Range("C8").Value = 13045.0932412616
$ws.Range("D8").Value = 12075.86
$ws.Range("E8").Value = 9107.12339392322
$ws.Range("F8").Value = 419.848193132701

$ws.Range("A9").Value = 46036
$ws.Range("B9").Value = 13378.423123507
$ws.Range("C9").Value = 12520.8328441972
$ws.Range("D9").Value = 12075.86
$ws.Range("E9").Value = 9107.12339392322
$ws.Range("F9").Value = 398.004009921682

$ws.Range("A10").Value = 46037
$ws.Range("B10").Value = 13378.423123507
$ws.Range("C10").Value = 12109.3353976805
$ws.Range("D10").Value = 12075.86
$ws.Range("E10").Value = 9107.12339392322
$ws.Range("F10").Value = 380.858282983486

$ws.Range("A11").Value = 46038
$ws.Range("B11").Value = 13378.423123507
$ws.Range("C11").Value = 11229.6404127551
$ws.Range("D11").Value = 12075.86
$ws.Range("E11").Value = 9107.08331334923
$ws.Range("F11").Value = 344.202655254347

$ws.Range("A12").Value = 46039
$ws.Range("B12").Value = 5471.00037786234
$ws.Range("C12").Value = 7767.17032312007
$ws.Range("D12").Value = 12075.86
$ws.Range("E12").Value = 8692.27239763912
$ws.Range("F12").Value = 182.649280031633

$ws.Range("A13").Value = 46040
$ws.Range("B13").Value = 5302.34398887746
$ws.Range("C13").Value = 7969.83077174213
$ws.Range("D13").Value = 12075.86
$ws.Range("E13").Value = 8679.8565018195
$ws.Range("F13").Value = 190.576136398402

$ws.Range("A14").Value = 46041
$ws.Range("B14").Value = 13071.6216038713
$ws.Range("C14").Value = 12207.3773095449
$ws.Range("D14").Value = 12075.86
$ws.Range("E14").Value = 8882.59744208774
$ws.Range("F14").Value = 375.58811465136

$ws.Range("A15").Value = 46042
$ws.Range("B15").Value = 13071.6216038713
$ws.Range("C15").Value = 12470.6014130971
$ws.Range("D15").Value = 12075.86
$ws.Range("E15").Value = 8882.59744208774
$ws.Range("F15").Value = 386.5557856327
